$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 27781278
$ws.Range("I62").Value = 37041036
$ws.Range("K62").Value = 37041036
$ws.Range("M62").Value = -37040412
# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 27781278
$ws.Range("I65").Value = 37041036
$ws.Range("K65").Value = 185205180
$ws.Range("M65").Value = -185202060
# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 4115.231
$ws.Range("I80").Value = 5550.4443
$ws.Range("J80").Value = 886
$ws.Range("K80").Value = 16651.3329
$ws.Range("L80").Value = 2658
$ws.Range("M80").Value = -15653.3329
$ws.Range("N80").Value = -4654
# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 4115.231
$ws.Range("I83").Value = 5550.4443
$ws.Range("J83").Value = 886
$ws.Range("K83").Value = 49953.9987
$ws.Range("L83").Value = 7974
$ws.Range("M83").Value = -44961.9987
$ws.Range("N83").Value = -17958
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 2585.279
$ws.Range("I98").Value = 2679
$ws.Range("J98").Value = 1671.5
$ws.Range("K98").Value = 2679
$ws.Range("L98").Value = 1671.5
$ws.Range("M98").Value = -1181
$ws.Range("N98").Value = -4667.5
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 2585.279
$ws.Range("I122").Value = 2679
$ws.Range("J122").Value = 1671.5
$ws.Range("K122").Value = 8037
$ws.Range("L122").Value = 5014.5
$ws.Range("M122").Value = -5587
$ws.Range("N122").Value = -9914.5
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 741.59576
$ws.Range("I129").Value = 484.45456
$ws.Range("J129").Value = 820.1667
$ws.Range("K129").Value = 1453.36368
$ws.Range("L129").Value = 2460.5001
$ws.Range("M129").Value = 3546.63632
$ws.Range("N129").Value = -12460.5001
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 7414675.5
$ws.Range("I132").Value = 12826677
$ws.Range("J132").Value = 8778.526
$ws.Range("K132").Value = 38480031
$ws.Range("L132").Value = 26335.578
$ws.Range("M132").Value = -38477501
$ws.Range("N132").Value = -31395.578
# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 40001690
$ws.Range("I135").Value = 482
$ws.Range("K135").Value = 4338
$ws.Range("M135").Value = -1803
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 544977.25
$ws.Range("I138").Value = 1111
$ws.Range("J138").Value = 806838.75
$ws.Range("K138").Value = 3333
$ws.Range("L138").Value = 2420516.25
$ws.Range("M138").Value = 1807
$ws.Range("N138").Value = -2430796.25

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5655.087
$ws.Range("I32").Value = 5230
$ws.Range("J32").Value = 15007
$ws.Range("K32").Value = 5230
$ws.Range("L32").Value = 15007
$ws.Range("M32").Value = -4943
$ws.Range("N32").Value = -15581
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1249.3077
$ws.Range("J45").Value = 1559.8
$ws.Range("L45").Value = 1559.8
$ws.Range("N45").Value = -2313.8
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 90910824
$ws.Range("I61").Value = 90910824
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 90910824
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -90910612
$ws.Range("N61").ClearContents()
# Row 92 (Leve Item ID 18050)
$ws.Range("H92").Value = 2507775
$ws.Range("J92").Value = 2507775
$ws.Range("L92").Value = 2507775
$ws.Range("N92").Value = -2512767
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 16669102
$ws.Range("I102").Value = 16669102
$ws.Range("K102").Value = 16669102
$ws.Range("M102").Value = -16667480
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 920.64703
$ws.Range("I122").Value = 935.7143
$ws.Range("J122").Value = 850.3333
$ws.Range("K122").Value = 2807.1429
$ws.Range("L122").Value = 2550.9999
$ws.Range("M122").Value = -357.1428999999998
$ws.Range("N122").Value = -7450.9999
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3263.1875
$ws.Range("I132").Value = 2654.5454
$ws.Range("K132").Value = 7963.6362
$ws.Range("M132").Value = -5433.6362
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 90910824
$ws.Range("I136").Value = 90910824
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 272732472
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -272729922
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 1225.4
$ws.Range("I107").Value = 1044.0625
$ws.Range("J107").Value = 1950.75
$ws.Range("K107").Value = 1044.0625
$ws.Range("L107").Value = 1950.75
$ws.Range("M107").Value = 875.9375
$ws.Range("N107").Value = -5790.75
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 6626
$ws.Range("I134").Value = 1095.0625
$ws.Range("K134").Value = 3285.1875
$ws.Range("M134").Value = -750.1875

$ws = $wb.Worksheets.Item("CRP")
# Row 118 (Leve Item ID 26112)
$ws.Range("H118").Value = 41950
$ws.Range("J118").Value = 41950
$ws.Range("L118").Value = 41950
$ws.Range("N118").Value = -45264
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 31253128
$ws.Range("I134").Value = 3982.2
$ws.Range("K134").Value = 11946.6
$ws.Range("M134").Value = -9411.599999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 411.9091
$ws.Range("I5").Value = 365.94446
$ws.Range("J5").Value = 618.75
$ws.Range("K5").Value = 1097.83338
$ws.Range("L5").Value = 1856.25
$ws.Range("M5").Value = -985.83338
$ws.Range("N5").Value = -2080.25
# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 266.22223
$ws.Range("I86").Value = 249.57143
$ws.Range("J86").Value = 324.5
$ws.Range("K86").Value = 748.71429
$ws.Range("L86").Value = 973.5
$ws.Range("M86").Value = 437.28571
$ws.Range("N86").Value = -3345.5
# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 266.22223
$ws.Range("I89").Value = 249.57143
$ws.Range("J89").Value = 324.5
$ws.Range("K89").Value = 2246.14287
$ws.Range("L89").Value = 2920.5
$ws.Range("M89").Value = 3681.85713
$ws.Range("N89").Value = -14776.5
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 16131931
$ws.Range("J131").Value = 3458.0784
$ws.Range("L131").Value = 10374.2352
$ws.Range("N131").Value = -20454.2352
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 411.9091
$ws.Range("I135").Value = 365.94446
$ws.Range("J135").Value = 618.75
$ws.Range("K135").Value = 3293.50014
$ws.Range("L135").Value = 5568.75
$ws.Range("M135").Value = -758.5001400000001
$ws.Range("N135").Value = -10638.75

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 28128452
$ws.Range("I70").Value = 25003520
$ws.Range("J70").Value = 33336670
$ws.Range("K70").Value = 25003520
$ws.Range("L70").Value = 33336670
$ws.Range("M70").Value = -25003250
$ws.Range("N70").Value = -33337210
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 28128452
$ws.Range("I73").Value = 25003520
$ws.Range("J73").Value = 33336670
$ws.Range("K73").Value = 25003520
$ws.Range("L73").Value = 33336670
$ws.Range("M73").Value = -25002584
$ws.Range("N73").Value = -33338542
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 644.1429000000001
$ws.Range("I97").Value = 627.25
$ws.Range("K97").Value = 627.25
$ws.Range("M97").Value = -131.25
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1019.4762
$ws.Range("I113").Value = 967.4167
$ws.Range("J113").Value = 1088.8889
$ws.Range("K113").Value = 967.4167
$ws.Range("L113").Value = 1088.8889
$ws.Range("M113").Value = 1202.5833
$ws.Range("N113").Value = -5428.8889
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3280.842
$ws.Range("I132").Value = 3083.625
$ws.Range("K132").Value = 9250.875
$ws.Range("M132").Value = -6720.875

$ws = $wb.Worksheets.Item("LTW")
# Row 44 (Leve Item ID 3658)
$ws.Range("H44").Value = 11266.667
$ws.Range("J44").Value = 11266.667
$ws.Range("L44").Value = 11266.667
$ws.Range("N44").Value = -12178.667
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 1094
$ws.Range("I61").Value = 1035.6428
$ws.Range("K61").Value = 1035.6428
$ws.Range("M61").Value = -833.6428000000001
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 1760.8
$ws.Range("I100").Value = 1301.3334
$ws.Range("J100").Value = 2450
$ws.Range("K100").Value = 1301.3334
$ws.Range("L100").Value = 2450
$ws.Range("M100").Value = -760.3334
$ws.Range("N100").Value = -3532
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 1094
$ws.Range("I113").Value = 1035.6428
$ws.Range("K113").Value = 1035.6428
$ws.Range("M113").Value = 1134.3572
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 62535376
$ws.Range("I122").Value = 83367230
$ws.Range("J122").Value = 39800
$ws.Range("K122").Value = 250101690
$ws.Range("L122").Value = 119400
$ws.Range("M122").Value = -250099240
$ws.Range("N122").Value = -124300
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 52512.5
$ws.Range("I132").Value = 14835.375
$ws.Range("J132").Value = 74042.28999999999
$ws.Range("K132").Value = 44506.125
$ws.Range("L132").Value = 222126.87
$ws.Range("M132").Value = -41976.125
$ws.Range("N132").Value = -227186.87

$ws = $wb.Worksheets.Item("WVR")
# Row 92 (Leve Item ID 18088)
$ws.Range("H92").Value = 20250
$ws.Range("J92").Value = 20250
$ws.Range("L92").Value = 20250
$ws.Range("N92").Value = -25242
# Row 99 (Leve Item ID 19640)
$ws.Range("H99").Value = 15500
$ws.Range("J99").Value = 15500
$ws.Range("L99").Value = 15500
$ws.Range("N99").Value = -21490
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 423.33334
$ws.Range("I113").Value = 275.55554
$ws.Range("J113").Value = 866.6667
$ws.Range("K113").Value = 826.66662
$ws.Range("L113").Value = 2600.0001
$ws.Range("M113").Value = 1343.33338
$ws.Range("N113").Value = -6940.0001
# Row 116 (Leve Item ID 26145)
$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 10418558
$ws.Range("I122").Value = 10871469
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 32614407
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -32611957
$ws.Range("N122").Value = -9700
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 47625932
$ws.Range("I126").Value = 58825750
$ws.Range("J126").Value = 26701
$ws.Range("K126").Value = 176477250
$ws.Range("L126").Value = 80103
$ws.Range("M126").Value = -176474780
$ws.Range("N126").Value = -85043
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2041.238
$ws.Range("I132").Value = 1861.421
$ws.Range("K132").Value = 5584.263
$ws.Range("M132").Value = -3054.263
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1048.0667
$ws.Range("I136").Value = 959.381
$ws.Range("K136").Value = 2878.143
$ws.Range("M136").Value = -328.143
